$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loadflow_Settings")

# Insert a new row at row 8 (shifts existing rows 8..59 down to 9..60),
# adding a new setting for "Automatic Tap Adjustment of Phase Shifters"
# (ldf.iPST_at) for backwards compatibility with existing Loadflow Settings.
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 2).Value = "Automatic Tap Adjustment of Phase Shifters"
$ws.Cells.Item(8, 3).Value = "ldf.iPST_at"
$ws.Cells.Item(8, 4).Value = 1

# Make Loadflow_Settings the active sheet/tab, with B11 selected.
$null = $ws.Select()
$null = $ws.Range("B11").Select()
